# Fruta / hortaliza, semanal
# Inserts two new data rows (new weekly observations for "Dina" variety,
# qualities "Especial" and "Primera") right before the existing row that
# used to be row 82, pushing all subsequent rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 82; Excel shifts rows 82:107 down to
# 84:109 and copies formatting (number formats, styles) from the row above
# into the freshly inserted blank rows, same as a native Excel "Insert
# Copied Cells"/"Insert Row" operation.
$ws.Rows("82:83").Insert()

# --- New row 82 ---
$ws.Cells.Item(82, 1).Value2 = 8
$ws.Cells.Item(82, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(82, 3).Value2 = "Coquimbo"
$ws.Cells.Item(82, 4).Value2 = 44932
$ws.Cells.Item(82, 5).Value2 = 4
$ws.Cells.Item(82, 6).Value2 = "Fruta"
$ws.Cells.Item(82, 7).Value2 = 100103
$ws.Cells.Item(82, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(82, 9).Value2 = 100103003
$ws.Cells.Item(82, 10).Value2 = "Damasco"
$ws.Cells.Item(82, 11).Value2 = "Dina"
$ws.Cells.Item(82, 12).Value2 = "Especial"
$ws.Cells.Item(82, 13).Value2 = 60
$ws.Cells.Item(82, 14).Value2 = 19000
$ws.Cells.Item(82, 15).Value2 = 20000
$ws.Cells.Item(82, 16).Value2 = 19500
$ws.Cells.Item(82, 17).Value2 = "`$/caja 16 kilos"
$ws.Cells.Item(82, 18).Value2 = "Región Metropolitana"
$ws.Cells.Item(82, 19).Value2 = 1219
$ws.Cells.Item(82, 20).Value2 = 16

# --- New row 83 ---
$ws.Cells.Item(83, 1).Value2 = 8
$ws.Cells.Item(83, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(83, 3).Value2 = "Coquimbo"
$ws.Cells.Item(83, 4).Value2 = 44932
$ws.Cells.Item(83, 5).Value2 = 4
$ws.Cells.Item(83, 6).Value2 = "Fruta"
$ws.Cells.Item(83, 7).Value2 = 100103
$ws.Cells.Item(83, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(83, 9).Value2 = 100103003
$ws.Cells.Item(83, 10).Value2 = "Damasco"
$ws.Cells.Item(83, 11).Value2 = "Dina"
$ws.Cells.Item(83, 12).Value2 = "Primera"
$ws.Cells.Item(83, 13).Value2 = 40
$ws.Cells.Item(83, 14).Value2 = 15000
$ws.Cells.Item(83, 15).Value2 = 16000
$ws.Cells.Item(83, 16).Value2 = 15500
$ws.Cells.Item(83, 17).Value2 = "`$/caja 16 kilos"
$ws.Cells.Item(83, 18).Value2 = "Región Metropolitana"
$ws.Cells.Item(83, 19).Value2 = 969
$ws.Cells.Item(83, 20).Value2 = 16
